{"js": "// Insert an \"Overview\" section after the Catchup bullet item that ends with\n// \"...ethical approval for user testing with the public\", and before the\n// \"Next Sprint Plans\" heading.\n//\n// New content (4 paragraphs, plain \"Normal\" style - i.e. not part of the\n// bulleted list):\n//   1) (blank paragraph)\n//   2) \"Overview\"  (bold)\n//   3) \"We caught up with what everyone has been doing, reviewed the\n//       project plan and Trello board and discussed what we must do next.\"\n//   4) (blank paragraph)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph by its (distinctive) text rather than by a\n// hard-coded index, so the script is resilient to minor structural drift.\nconst anchorText = \"Kieran is also working together with the other teams to get ethical approval for user testing with the public\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Could not locate anchor paragraph for Overview insertion\");\n}\n\n// 1) Blank paragraph right after the bullet list item. Force the style to\n//    \"Normal\" so it leaves the bulleted (\"ListParagraph\") list instead of\n//    inheriting the numbering/indent of the anchor paragraph.\nconst blank1 = anchor.insertParagraph(\"\", Word.InsertLocation.after);\nblank1.style = \"Normal\";\nawait context.sync();\n\n// 2) \"Overview\" heading-like paragraph, bold text.\nconst overviewPara = blank1.insertParagraph(\"Overview\", Word.InsertLocation.after);\noverviewPara.style = \"Normal\";\noverviewPara.font.bold = true;\n// Also flag the complex-script bold flag (w:bCs) to mirror Word's usual\n// pairing of w:b/w:bCs when bolding text from the UI.\noverviewPara.font.boldBidirectional = true;\nawait context.sync();\n\n// 3) Body paragraph summarizing the catch-up.\nconst summaryPara = overviewPara.insertParagraph(\n  \"We caught up with what everyone has been doing, reviewed the project plan and Trello board and discussed what we must do next.\",\n  Word.InsertLocation.after\n);\nsummaryPara.style = \"Normal\";\nawait context.sync();\n\n// 4) Trailing blank paragraph separating this section from \"Next Sprint Plans\".\nconst blank2 = summaryPara.insertParagraph(\"\", Word.InsertLocation.after);\nblank2.style = \"Normal\";\nawait context.sync();\n", "ps1": "# Insert an \"Overview\" section after the Catchup bullet item that ends with\n# \"...ethical approval for user testing with the public\", and before the\n# \"Next Sprint Plans\" heading.\n#\n# New content (4 paragraphs, plain \"Normal\" style - i.e. not part of the\n# bulleted list):\n#   1) (blank paragraph)\n#   2) \"Overview\"  (bold)\n#   3) \"We caught up with what everyone has been doing, reviewed the\n#       project plan and Trello board and discussed what we must do next.\"\n#   4) (blank paragraph)\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph by its (distinctive) text rather than a\n# hard-coded index, so the script is resilient to minor structural drift.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*ethical approval for user testing with the public*\") {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not locate anchor paragraph for Overview insertion\"\n}\n\n$anchorPara = $d.Paragraphs.Item($anchorIndex)\n\n# 1) Blank paragraph right after the bullet list item. Remove the inherited\n#    numbering and force the style to \"Normal\" so it leaves the bulleted\n#    (\"List Paragraph\") list instead of staying part of it.\n$anchorPara.Range.InsertParagraphAfter()\n$blank1 = $d.Paragraphs.Item($anchorIndex + 1)\n$blank1.Range.ListFormat.RemoveNumbers()\n$blank1.Style = \"Normal\"\n\n# 2) \"Overview\" paragraph (bold). Insert its trailing paragraph break FIRST\n#    so the paragraph mark that carries over to paragraph 3 below does not\n#    pick up the bold formatting we apply afterwards.\n$blank1.Range.InsertParagraphAfter()\n$overviewPara = $d.Paragraphs.Item($anchorIndex + 2)\n$overviewPara.Range.ListFormat.RemoveNumbers()\n$overviewPara.Style = \"Normal\"\n\n$overviewPara.Range.InsertParagraphAfter()\n$summaryPara = $d.Paragraphs.Item($anchorIndex + 3)\n$summaryPara.Range.ListFormat.RemoveNumbers()\n$summaryPara.Style = \"Normal\"\n\n$overviewPara.Range.Text = \"Overview\"\n$overviewPara.Range.Bold = 1\n# Also flag the complex-script bold flag (w:bCs) to mirror Word's usual\n# pairing of w:b/w:bCs when bolding text from the UI.\n$overviewPara.Range.Font.BoldBi = 1\n\n# 3) Body paragraph summarizing the catch-up (not bold).\n$summaryPara.Range.Text = \"We caught up with what everyone has been doing, reviewed the project plan and Trello board and discussed what we must do next.\"\n\n# 4) Trailing blank paragraph separating this section from \"Next Sprint Plans\".\n$summaryPara.Range.InsertParagraphAfter()\n$blank2 = $d.Paragraphs.Item($anchorIndex + 4)\n$blank2.Range.ListFormat.RemoveNumbers()\n$blank2.Style = \"Normal\"\n"}
